# Update "想去人数" (people interested) counts across the relevant sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 3450
$wsExpo.Range("F5").Value = 6987
$wsExpo.Range("F6").Value = 2505
$wsExpo.Range("F11").Value = 81
$wsExpo.Range("F12").Value = 36

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 24

# Sheet "全部类型" (All types - aggregate of the above)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3450
$wsAll.Range("F3").Value = 24
$wsAll.Range("F6").Value = 6987
$wsAll.Range("F7").Value = 2505
$wsAll.Range("F12").Value = 81
$wsAll.Range("F13").Value = 36
